# Insert a new "TECHNICAL DETAILS" table (Capture/Detection Antibodies,
# Specificity, Standard Protein, Cross-reactivity vs. technical_details_table
# Jinja placeholders) right after the "TECHNICAL DETAILS" heading paragraph
# and before the existing "{{ technical_details }}" paragraph.

$d = $word.ActiveDocument

# Find the "TECHNICAL DETAILS" heading paragraph.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.Trim() -eq "TECHNICAL DETAILS") {
        $target = $p
        break
    }
}

if ($target -eq $null) {
    throw "Could not locate the TECHNICAL DETAILS heading paragraph"
}

# Build a fresh, collapsed Range sitting right at the end of that heading
# paragraph (i.e. right before the following "{{ technical_details }}"
# paragraph) and insert the table's WordOpenXML fragment there.
$insertPos = $target.Range.End
$rng = $d.Range($insertPos, $insertPos)

$xml = @'
<w:tbl xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:tblPr>
    <w:tblStyle w:val="TableGrid"/>
    <w:tblW w:type="auto" w:w="0"/>
    <w:tblLook w:firstColumn="1" w:firstRow="1" w:lastColumn="0" w:lastRow="0" w:noHBand="0" w:noVBand="1" w:val="04A0"/>
  </w:tblPr>
  <w:tblGrid>
    <w:gridCol w:w="5400"/>
    <w:gridCol w:w="5400"/>
  </w:tblGrid>
  <w:tr>
    <w:tc>
      <w:tcPr>
        <w:tcW w:type="dxa" w:w="4320"/>
      </w:tcPr>
      <w:p>
        <w:r>
          <w:t>Capture/Detection Antibodies</w:t>
        </w:r>
      </w:p>
    </w:tc>
    <w:tc>
      <w:tcPr>
        <w:tcW w:type="dxa" w:w="4320"/>
      </w:tcPr>
      <w:p>
        <w:r>
          <w:t>{{ technical_details_table[0].value }}</w:t>
        </w:r>
      </w:p>
    </w:tc>
  </w:tr>
  <w:tr>
    <w:tc>
      <w:tcPr>
        <w:tcW w:type="dxa" w:w="4320"/>
      </w:tcPr>
      <w:p>
        <w:r>
          <w:t>Specificity</w:t>
        </w:r>
      </w:p>
    </w:tc>
    <w:tc>
      <w:tcPr>
        <w:tcW w:type="dxa" w:w="4320"/>
      </w:tcPr>
      <w:p>
        <w:r>
          <w:t>{{ technical_details_table[1].value }}</w:t>
        </w:r>
      </w:p>
    </w:tc>
  </w:tr>
  <w:tr>
    <w:tc>
      <w:tcPr>
        <w:tcW w:type="dxa" w:w="4320"/>
      </w:tcPr>
      <w:p>
        <w:r>
          <w:t>Standard Protein</w:t>
        </w:r>
      </w:p>
    </w:tc>
    <w:tc>
      <w:tcPr>
        <w:tcW w:type="dxa" w:w="4320"/>
      </w:tcPr>
      <w:p>
        <w:r>
          <w:t>{{ technical_details_table[2].value }}</w:t>
        </w:r>
      </w:p>
    </w:tc>
  </w:tr>
  <w:tr>
    <w:tc>
      <w:tcPr>
        <w:tcW w:type="dxa" w:w="4320"/>
      </w:tcPr>
      <w:p>
        <w:r>
          <w:t>Cross-reactivity</w:t>
        </w:r>
      </w:p>
    </w:tc>
    <w:tc>
      <w:tcPr>
        <w:tcW w:type="dxa" w:w="4320"/>
      </w:tcPr>
      <w:p>
        <w:r>
          <w:t>{{ technical_details_table[3].value }}</w:t>
        </w:r>
      </w:p>
    </w:tc>
  </w:tr>
</w:tbl>
'@

$rng.InsertXML($xml)

Write-Output "Inserted TECHNICAL DETAILS table"
